$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 51999.5
$ws.Range("J3").Value = 51999.5
$ws.Range("L3").Value = 51999.5
$ws.Range("N3").Value = -52227.5

$ws.Range("H9").Value = 584.625
$ws.Range("I9").Value = 603.8570999999999
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 603.8570999999999
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = -434.8570999999999
$ws.Range("N9").Value = -788

$ws.Range("H18").Value = 15424.75
$ws.Range("I18").Value = 8899.666999999999
$ws.Range("K18").Value = 8899.666999999999
$ws.Range("M18").Value = -8615.666999999999

$ws.Range("H33").Value = 115.666664
$ws.Range("I33").Value = 108.8
$ws.Range("K33").Value = 108.8
$ws.Range("M33").Value = 120.2

$ws.Range("H40").Value = 4910.778
$ws.Range("J40").Value = 2991.5833
$ws.Range("L40").Value = 2991.5833
$ws.Range("N40").Value = -3341.5833

$ws.Range("H92").Value = 1288.75
$ws.Range("I92").Value = 462.69232
$ws.Range("J92").Value = 4868.3335
$ws.Range("K92").Value = 462.69232
$ws.Range("L92").Value = 4868.3335
$ws.Range("M92").Value = 785.30768
$ws.Range("N92").Value = -7364.3335

$ws.Range("H102").Value = 51999.5
$ws.Range("J102").Value = 51999.5
$ws.Range("L102").Value = 51999.5
$ws.Range("N102").Value = -58489.5

$ws.Range("H129").Value = 3341.1667
$ws.Range("I129").Value = 3136.75
$ws.Range("K129").Value = 9410.25
$ws.Range("M129").Value = -4410.25

$ws.Range("H132").Value = 463544.7
$ws.Range("I132").Value = 546273.9399999999
$ws.Range("J132").Value = 8534
$ws.Range("K132").Value = 1638821.82
$ws.Range("L132").Value = 25602
$ws.Range("M132").Value = -1636291.82
$ws.Range("N132").Value = -30662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H32").Value = 1824686.9
$ws.Range("I32").Value = 2118242.2
$ws.Range("K32").Value = 2118242.2
$ws.Range("M32").Value = -2117955.2

$ws.Range("H45").Value = 2395.926
$ws.Range("I45").Value = 2689.0527
$ws.Range("K45").Value = 2689.0527
$ws.Range("M45").Value = -2312.0527

$ws.Range("H61").Value = 7600.048
$ws.Range("I61").Value = 4043.2856
$ws.Range("K61").Value = 4043.2856
$ws.Range("M61").Value = -3831.2856

$ws.Range("H110").Value = 2014.0883
$ws.Range("I110").Value = 2000.6207
$ws.Range("K110").Value = 2000.6207
$ws.Range("M110").Value = 44.37930000000006

$ws.Range("H132").Value = 373921.12
$ws.Range("I132").Value = 508298.8
$ws.Range("K132").Value = 1524896.4
$ws.Range("M132").Value = -1522366.4

$ws.Range("H136").Value = 7600.048
$ws.Range("I136").Value = 4043.2856
$ws.Range("K136").Value = 12129.8568
$ws.Range("M136").Value = -9579.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 7399.2666
$ws.Range("J8").Value = 7213.5
$ws.Range("L8").Value = 7213.5
$ws.Range("N8").Value = -7493.5

$ws.Range("H20").Value = 2067.7932
$ws.Range("I20").Value = 2342.889
$ws.Range("J20").Value = 1617.6364
$ws.Range("K20").Value = 2342.889
$ws.Range("L20").Value = 1617.6364
$ws.Range("M20").Value = -2095.889
$ws.Range("N20").Value = -2111.6364

$ws.Range("H134").Value = 646500.5600000001
$ws.Range("I134").Value = 770142.9
$ws.Range("K134").Value = 2310428.7
$ws.Range("M134").Value = -2307893.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11057.75
$ws.Range("I31").Value = 4274.2856
$ws.Range("J31").Value = 16333.777
$ws.Range("K31").Value = 4274.2856
$ws.Range("L31").Value = 16333.777
$ws.Range("M31").Value = -3979.2856
$ws.Range("N31").Value = -16923.777

$ws.Range("H34").Value = 11057.75
$ws.Range("I34").Value = 4274.2856
$ws.Range("J34").Value = 16333.777
$ws.Range("K34").Value = 4274.2856
$ws.Range("L34").Value = 16333.777
$ws.Range("M34").Value = -4072.2856
$ws.Range("N34").Value = -16737.777

$ws.Range("H58").Value = 777179.0600000001
$ws.Range("J58").Value = 3820.8333
$ws.Range("L58").Value = 3820.8333
$ws.Range("N58").Value = -4226.8333

$ws.Range("H94").Value = 1288.2858
$ws.Range("I94").Value = 778.7143
$ws.Range("J94").Value = 1797.8572
$ws.Range("K94").Value = 778.7143
$ws.Range("L94").Value = 1797.8572
$ws.Range("M94").Value = -327.7143
$ws.Range("N94").Value = -2699.8572

$ws.Range("H134").Value = 33756.19
$ws.Range("I134").Value = 44062.5
$ws.Range("K134").Value = 132187.5
$ws.Range("M134").Value = -129652.5

$ws.Range("H136").Value = 777179.0600000001
$ws.Range("J136").Value = 3820.8333
$ws.Range("L136").Value = 11462.4999
$ws.Range("N136").Value = -16562.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5364.4287
$ws.Range("J3").Value = 7874
$ws.Range("L3").Value = 23622
$ws.Range("N3").Value = -23846

$ws.Range("H105").Value = 15303.1
$ws.Range("J105").Value = 16336.777
$ws.Range("L105").Value = 49010.331
$ws.Range("N105").Value = -54252.331

$ws.Range("H113").Value = 3084.3547
$ws.Range("I113").Value = 1858.25
$ws.Range("J113").Value = 3510.8262
$ws.Range("K113").Value = 5574.75
$ws.Range("L113").Value = 10532.4786
$ws.Range("M113").Value = -3404.75
$ws.Range("N113").Value = -14872.4786

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1406.3529
$ws.Range("I113").Value = 1406.3529
$ws.Range("K113").Value = 1406.3529
$ws.Range("M113").Value = 763.6470999999999

$ws.Range("H126").Value = 1193484.6
$ws.Range("I126").Value = 1853812.9
$ws.Range("K126").Value = 5561438.699999999
$ws.Range("M126").Value = -5558968.699999999

$ws.Range("H132").Value = 2479.76
$ws.Range("I132").Value = 2329.3333
$ws.Range("K132").Value = 6987.999899999999
$ws.Range("M132").Value = -4457.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 70250
$ws.Range("J22").Value = 8266.666999999999
$ws.Range("L22").Value = 8266.666999999999
$ws.Range("N22").Value = -8856.666999999999

$ws.Range("H27").Value = 70250
$ws.Range("J27").Value = 8266.666999999999
$ws.Range("L27").Value = 8266.666999999999
$ws.Range("N27").Value = -8480.666999999999

$ws.Range("H54").Value = 30000
$ws.Range("J54").Value = 30000
$ws.Range("L54").Value = 30000
$ws.Range("N54").Value = -31288

$ws.Range("H55").Value = 1812.0769
$ws.Range("I55").Value = 733
$ws.Range("K55").Value = 733
$ws.Range("M55").Value = -560

$ws.Range("H122").Value = 49565.914
$ws.Range("I122").Value = 4447.9414
$ws.Range("J122").Value = 177400.17
$ws.Range("K122").Value = 13343.8242
$ws.Range("L122").Value = 532200.51
$ws.Range("M122").Value = -10893.8242
$ws.Range("N122").Value = -537100.51

$ws.Range("H132").Value = 964845.4399999999
$ws.Range("I132").Value = 1444221.8
$ws.Range("K132").Value = 4332665.4
$ws.Range("M132").Value = -4330135.4

$ws.Range("H136").Value = 5481.4
$ws.Range("I136").Value = 5301.75
$ws.Range("J136").Value = 6200
$ws.Range("K136").Value = 15905.25
$ws.Range("L136").Value = 18600
$ws.Range("M136").Value = -13355.25
$ws.Range("N136").Value = -23700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 800
$ws.Range("I7").Value = 800
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = -687
$ws.Range("M7").ClearContents()

$ws.Range("H9").Value = 1007
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1007
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = 1007
$ws.Range("N9").Value = -1287
$ws.Range("L9").ClearContents()

$ws.Range("H107").Value = 1227.7
$ws.Range("I107").Value = 1077.3043
$ws.Range("K107").Value = 3231.9129
$ws.Range("M107").Value = -1311.9129

$ws.Range("H136").Value = 11177373
$ws.Range("I136").Value = 19998938
$ws.Range("J136").Value = 3390.0667
$ws.Range("K136").Value = 59996814
$ws.Range("L136").Value = 10170.2001
$ws.Range("M136").Value = -59994264
$ws.Range("N136").Value = -15270.2001
